# Adds the PTC_FR dataset (as a new transposed "record" in column N, one
# value per header row 1..14) to the test_datasets sheet, and extends the
# row-index column A down to match (values 1..14), mirroring the existing
# MUTAG row/column layout. Matches commit "Added NX-GED; Cleand Calculators;".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- N1: new header cell for the extra record column, same header style as
#     B1:M1 (bold / bordered / centered), value 0 -------------------------
$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Cells.Item(1, 14).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(1, 14).Value = 0

# --- A3:A16: extend the row-index column, same style as A2 (s="1") -------
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Range($ws.Cells.Item(3, 1), $ws.Cells.Item(16, 1)).PasteSpecial(-4122) | Out-Null

for ($row = 3; $row -le 16; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 2
}

# --- N3:N16: the PTC_FR dataset values (one per header field) ------------
$ptcFrValues = @(
    'PTC_FR',
    'TUD',
    'Bioinformatics',
    351,
    2,
    $true,
    $true,
    14.55840455840456,
    15.002849002849,
    '{0: 0.6552706552706553, 1: 0.34472934472934474}',
    '[''0'', ''1'', ''10'', ''11'', ''12'', ''13'', ''14'', ''15'', ''16'', ''17'', ''18'', ''2'', ''3'', ''4'', ''5'', ''6'', ''7'', ''8'', ''9'']',
    '[''0'', ''1'', ''2'', ''3'']',
    "[   1   23    1    1    2    1    3    1    1    1    1  721  408   29`n 3493  269  100   44   10]",
    '[   7  436 2772 2051]'
)

for ($i = 0; $i -lt $ptcFrValues.Count; $i++) {
    $ws.Cells.Item($i + 3, 14).Value = $ptcFrValues[$i]
}

# N15's value contains an embedded newline, which triggers Excel's
# auto row-height ("wrap to fit") adjustment; restore the default height.
$ws.Rows.Item(15).AutoFit() | Out-Null
